# Add arrow-label columns (A_arrow, B_arrow, C_arrow) to the "axes" sheet,
# shifting the existing "Title" column from D to G, widen the new columns,
# and make "axes" the active/selected sheet.

$wb = $excel.ActiveWorkbook
$axes = $wb.Worksheets.Item("axes")

# Remember the existing "Title" header/value (currently in column D) before
# the new columns push it over to column G.
$titleHeader = $axes.Range("D1").Value()
$titleValue  = $axes.Range("D2").Value()

# New header row for the three inserted columns.
$axes.Range("D1").Value = "A_arrow"
$axes.Range("E1").Value = "B_arrow"
$axes.Range("F1").Value = "C_arrow"
$axes.Range("G1").Value = $titleHeader

# New data row: mirror the A/B/C axis values (An, Ab, Or), then restore the
# Title text in its new column G.
$axes.Range("D2").Value = $axes.Range("A2").Value()
$axes.Range("E2").Value = $axes.Range("B2").Value()
$axes.Range("F2").Value = $axes.Range("C2").Value()
$axes.Range("G2").Value = $titleValue

# Widen the three new columns.
$axes.Columns.Item(4).ColumnWidth = 13.3
$axes.Columns.Item(5).ColumnWidth = 13.3
$axes.Columns.Item(6).ColumnWidth = 13.3

$axes.Range("F2").Select() | Out-Null

# "axes" becomes the active/selected tab.
$axes.Activate()
